# rd3stats.xlsx: expanded tree data table (#21)
#
# - packages: bump version/date string in the description cell
# - entities: replace the "json" entity row with the new "treedata" entity
# - attributes: re-point idAttribute/labelAttribute/lookupAttribute columns,
#   add a new "nillable" column, and add a new "dataType" column (H) holding
#   what used to live in the old E (dataType) column; replace the two
#   "rd3stats_json" attribute rows with three "rd3stats_treedata" rows
#   (subjectID, familyID, json)

$wb = $excel.ActiveWorkbook

# --- Sheet "packages" ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("packages")
$ws1.Range("C2").Value = "Additional summaries and processed data (v1.1.0, 2022-06-27)"

# --- Sheet "entities" ----------------------------------------------------
$ws2 = $wb.Worksheets.Item("entities")
$ws2.Range("B2").Value = "treedata"
$ws2.Range("C2").Value = "Patient Tree Data"
$ws2.Range("D2").Value = "JSON stringified objects of sample-experiment links per subject"

# --- Sheet "attributes" ---------------------------------------------------
$ws3 = $wb.Worksheets.Item("attributes")

# Insert a new column H (inherits formatting, e.g. the bordered header
# style, from the adjacent column G automatically).
$ws3.Columns.Item(8).Insert()

# Header row: E/F shift meaning, G becomes "nillable", H becomes "dataType"
$ws3.Range("E1").Value = "labelAttribute"
$ws3.Range("F1").Value = "lookupAttribute"
$ws3.Range("G1").Value = "nillable"
$ws3.Range("H1").Value = "dataType"

# Row 2: subjectID
$ws3.Range("A2").Value = "rd3stats_treedata"
$ws3.Range("B2").Value = "subjectID"
$ws3.Range("C2").Value = "An individual who is the subject of personal data, persons to whom data refers, and from whom data are collected, processed, and stored."
$ws3.Range("D2").Value = $true
$ws3.Range("E2").Value = $true
$ws3.Range("F2").Value = $true
$ws3.Range("G2").Value = $false
$ws3.Range("H2").Value = "string"

# Row 3: familyID
$ws3.Range("A3").Value = "rd3stats_treedata"
$ws3.Range("B3").Value = "familyID"
$ws3.Range("C3").Value = "A domestic group, or a number of domestic groups linked through descent (demonstrated or stipulated) from a common ancestor, marriage, or adoption."
$ws3.Range("D3").Value = $false
$ws3.Range("E3").Value = $false
$ws3.Range("F3").Value = $false
$ws3.Range("G3").ClearContents()
$ws3.Range("H3").Value = "string"

# Row 4 (new): json
$ws3.Range("A4").Value = "rd3stats_treedata"
$ws3.Range("B4").Value = "json"
$ws3.Range("C4").Value = "json stringified object containing sample-experiment links"
$ws3.Range("D4").Value = $false
$ws3.Range("E4").Value = $false
$ws3.Range("F4").Value = $false
$ws3.Range("H4").Value = "text"
